$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-25 08:49:02"

$wsZhCn.Range("H4").Value = "2016-08-25 08:48:56"
$wsZhCn.Range("K4").Value = "2016-08-25 08:49:31"

$wsDeDe.Range("H4").Value = "2016-08-25 08:49:02"
$wsDeDe.Range("K4").Value = "2016-08-25 08:49:38"
